$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes on row 4 ---
$ws.Range("C4").Value = 7
$ws.Range("E4").Value = 9

# --- Strip the now-unused fill/border style from G4:G8 (back to default "Standard" style) ---
$ws.Range("G4:G8").Style = "Standard"

# --- Update the active selection on the sheet ---
$ws.Range("C5").Select()
